$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Agosto de 2020 a las 13:59"

# --- Update country stat rows (columns: B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5798983
$ws.Range("C4").Value = 2256
$ws.Range("E4").Value = 2492325
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 179240

# Iran (row 14)
$ws.Range("B14").Value = 356792
$ws.Range("C14").Value = 2028
$ws.Range("D14").Value = 307702
$ws.Range("E14").Value = 28588
$ws.Range("G14").Value = 126
$ws.Range("H14").Value = 20502

# Emiratos Arabes Unidos (row 45)
$ws.Range("B45").Value = 66617
$ws.Range("C45").Value = 424
$ws.Range("D45").Value = 58408
$ws.Range("E45").Value = 7837
$ws.Range("G45").Value = 2
$ws.Range("H45").Value = 372

# Estado de Palestina (row 76)
$ws.Range("B76").Value = 18476
$ws.Range("C76").Value = 163
$ws.Range("E76").Value = 7252

# Finlandia (row 100)
$ws.Range("B100").Value = 7906
$ws.Range("C100").Value = 35
$ws.Range("E100").Value = 472

# Eslovenia's case counts rose enough to overtake Tunez and Lituania in the
# ranking (table is sorted descending by Casos totales), so it moves from
# row 130 up to row 128; Tunez and Lituania each shift down one row with
# their data otherwise unchanged.
$ws.Range("A128").Value = "Eslovenia"
$ws.Range("B128").Value = 2617
$ws.Range("C128").Value = 43
$ws.Range("D128").Value = 2079
$ws.Range("E128").Value = 407
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 131

$ws.Range("A129").Value = "Tunez"
$ws.Range("B129").Value = 2607
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 1420
$ws.Range("E129").Value = 1123
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 64

$ws.Range("A130").Value = "Lituania"
$ws.Range("B130").Value = 2594
$ws.Range("C130").Value = 30
$ws.Range("D130").Value = 1766
$ws.Range("E130").Value = 744
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 1
$ws.Range("H130").Value = 84

# Vietnam (row 160)
$ws.Range("B160").Value = 1014
$ws.Range("C160").Value = 5
$ws.Range("E160").Value = 444

# Gibraltar (row 185)
$ws.Range("B185").Value = 239
$ws.Range("C185").Value = 8
$ws.Range("D185").Value = 201
$ws.Range("E185").Value = 38
